{"js": "// Remove the stale \"Ver no Jupiter / Salvar ...\" line, the copyright/footer\n// line, and the blank paragraph that separated them from the preceding\n// \"Requisitos\" text, as part of a routine site rebuild that dropped those\n// three trailing paragraphs right after the \"LOM3013: ...\" requirement line.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the anchor paragraph (\"LOM3013: Ci\u00eancia dos Materiais (Requisito\n// fraco)\") so the removal is tied to content rather than a brittle index.\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"LOM3013: Ci\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error('Could not find the \"LOM3013\" requirement paragraph.');\n}\n\n// The three paragraphs immediately following the anchor are: a blank\n// paragraph, the \"Ver no Jupiter ...\" line, and the \"\u00a9 2020 ...\" line.\n// Delete them (back to front so indices stay valid while deleting).\nconst toRemoveTexts = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\\u00A9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\",\n];\n\nconst candidates = [];\nfor (let offset = 1; offset <= 3 && anchorIndex + offset < items.length; offset++) {\n  candidates.push(anchorIndex + offset);\n}\n\nfor (let i = candidates.length - 1; i >= 0; i--) {\n  const idx = candidates[i];\n  const text = items[idx].text;\n  if (text === \"\" || toRemoveTexts.indexOf(text) !== -1) {\n    items[idx].delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the stale \"Ver no Jupiter / Salvar ...\" line, the copyright/footer\n# line, and the blank paragraph that separated them from the preceding\n# \"Requisitos\" text -- three trailing paragraphs that followed the\n# \"LOM3013: ...\" requirement line were dropped in this site rebuild.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph (\"LOM3013: Ci\u00eancia dos Materiais (Requisito\n# fraco)\") by content so the removal isn't tied to a brittle fixed index.\n$count = $d.Paragraphs.Count\n$anchorIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*LOM3013: Ci*\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not find the 'LOM3013' requirement paragraph.\"\n}\n\n# The paragraphs immediately after the anchor are: a blank paragraph, the\n# \"Ver no Jupiter ...\" line, and the \"\u00a9 2020 ...\" line. Identify them by\n# text, then delete back-to-front so indices stay valid while deleting.\n$removableTexts = @(\n    \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n    [char]0x00A9 + \" 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n)\n\n$candidates = New-Object System.Collections.ArrayList\nfor ($offset = 1; $offset -le 3; $offset++) {\n    $idx = $anchorIndex + $offset\n    if ($idx -le $d.Paragraphs.Count) {\n        [void]$candidates.Add($idx)\n    }\n}\n\nfor ($j = $candidates.Count - 1; $j -ge 0; $j--) {\n    $idx = $candidates[$j]\n    $para = $d.Paragraphs.Item($idx)\n    $text = $para.Range.Text.TrimEnd([char]13, [char]7)\n    if ($text -eq \"\" -or $removableTexts -contains $text) {\n        $para.Range.Delete()\n    }\n}\n"}
